# Applies the "Updated symbol list" commit: refreshes the Price (column D)
# figures for most rows, rewrites a couple of Volume(1h) labels (column E),
# and rotates the BKEXToken / CEJI / KickToken rows (41-43) into their new
# order with updated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while keeping it a genuine text value
# (the sheet stores every data cell as a string, even ones that look like
# numbers), and without leaving behind any stray per-cell style/format.
function Set-TextValue {
    param(
        $Sheet,
        [string]$Address,
        [string]$Value
    )
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

# ---- Column D (Price) refreshes ----
Set-TextValue $ws "D2"  "245.76"
Set-TextValue $ws "D3"  "21.99"
Set-TextValue $ws "D4"  "5.377"
Set-TextValue $ws "D5"  "0.05859"
Set-TextValue $ws "D6"  "3.388"
Set-TextValue $ws "D7"  "6.365"
Set-TextValue $ws "D8"  "0.8134"
Set-TextValue $ws "D9"  "1.002"
Set-TextValue $ws "D10" "0.1416"
Set-TextValue $ws "D11" "0.03880"
Set-TextValue $ws "D12" "0.07410"
Set-TextValue $ws "D13" "0.03043"
Set-TextValue $ws "D14" "4.159"
Set-TextValue $ws "D15" "0.09391"
Set-TextValue $ws "D16" "0.001585"
Set-TextValue $ws "D17" "0.04811"
Set-TextValue $ws "D18" "0.0005886"

# ---- Column E (Volume(1h)) label tweak ----
Set-TextValue $ws "E18" "17OneONEWorstin24h"

Set-TextValue $ws "D19" "0.006033"
Set-TextValue $ws "D20" "0.004092"
Set-TextValue $ws "D21" "0.0009886"
Set-TextValue $ws "D22" "0.0001498"
Set-TextValue $ws "D23" "3.697"
Set-TextValue $ws "D24" "2.225"
Set-TextValue $ws "D27" "0.0002492"
Set-TextValue $ws "D40" "0.03866"

# ---- Rows 41-43: BKEXToken / CEJI / KickToken rotate positions ----
# New row 41 = KickToken (was row 43's coin, with refreshed data)
Set-TextValue $ws "B41" "KickToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006412"
Set-TextValue $ws "E41" "40KickTokenKICK"

# New row 42 = BKEXToken (was row 41's coin, with refreshed data)
Set-TextValue $ws "B42" "BKEXToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1074"
Set-TextValue $ws "E42" "41BKEXTokenBKK"

# New row 43 = CEJI (was row 42's coin, with refreshed data)
Set-TextValue $ws "B43" "CEJI"
Set-TextValue $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002997"
Set-TextValue $ws "E43" "42CEJICEJI"

# ---- Remaining column D (Price) refreshes ----
Set-TextValue $ws "D44" "0.006677"
Set-TextValue $ws "D45" "0.00005619"
Set-TextValue $ws "D46" "0.00000000749"
Set-TextValue $ws "D47" "0.6995"
Set-TextValue $ws "D49" "0.00002098"
Set-TextValue $ws "D50" "0.01009"
